$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.170.37"
$ws.Range("E2").Value = "  -2.61%  "
$ws.Range("D3").Value = "3.135.53"
$ws.Range("E3").Value = "  -4.23%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.19%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.133.95"
$ws.Range("E8").Value = "  -4.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.443"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.37%  "
$ws.Range("E11").Value = "  -8.34%  "
$ws.Range("E12").Value = "  -6.45%  "
$ws.Range("D13").Value = "3.671.56"
$ws.Range("E13").Value = "  -4.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.128"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.43%  "
$ws.Range("D16").Value = "3.133.37"
$ws.Range("E16").Value = "  -4.27%  "
$ws.Range("D17").Value = "58.155.52"
$ws.Range("E17").Value = "  -2.82%  "
$ws.Range("E18").Value = "  -6.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "343.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.62%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.507"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.90%  "
$ws.Range("D26").Value = "3.259.13"
$ws.Range("E26").Value = "  -4.44%  "
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("E30").Value = "  -3.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  -7.55%  "
$ws.Range("E33").Value = "  -7.38%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "21.34"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.81"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.73%  "
$ws.Range("E38").Value = "  -5.61%  "
$ws.Range("E39").Value = "  -9.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0684"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.41%  "
$ws.Range("D41").Value = "3.165.95"
$ws.Range("E41").Value = "  -4.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.57%  "
$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("E45").Value = "  -7.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.82%  "
$ws.Range("E50").Value = "  -2.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.59%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0951"
$ws.Range("E29").Value = "  -6.15%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.52%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.277.25"
$ws.Range("E49").Value = "  -1.62%  "
